# Quarterly Performance.xlsx - remove the "Tenure Years" column (D) from the
# header row, shifting the remaining headers (Productivity % ... HR Comments)
# one column to the left, without disturbing the stored column width
# definitions (<cols>) that Excel's EntireColumn/Range "shift-left" delete
# would otherwise renumber.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift D1:K1 one cell to the left (i.e. drop the old D1 "Tenure Years"
# value, pulling every following header left by one column) by copying
# values directly cell-by-cell instead of using Range.Delete, which would
# also re-map the <cols> width metadata.
for ($col = 4; $col -le 10; $col++) {
    $ws.Cells.Item(1, $col).Value2 = $ws.Cells.Item(1, $col + 1).Value2
}

# The last column (K1) is now a duplicate of J1; clear it so the used range
# shrinks back down to column J.
$ws.Cells.Item(1, 11).ClearContents()

# Match the saved selection state: active cell D1, selection D1:J1.
$ws.Range("D1:J1").Select()
